# Applies the commit's change: inserts three new "general_college_subjects"
# columns (history, electives, cs) before the existing "arts" column, which
# shifts every column from the old R onward three places to the right
# (old R:AE -> new U:AH), and lower-cases the "Unknown" placeholder text
# for the per-school importance columns (D:J) while leaving the
# state.State column (K) as "Unknown".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new whole columns at R, shifting R:AE (and everything in
#    between) three columns to the right -> U:AH. xlShiftToRight = -4161.
$ws.Range("R1:T1048576").Insert(-4161)

# 2. Populate the headers for the three newly inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# 3. Populate the matching data-row cells (numeric, like the other
#    general_college_subjects columns) for the sample row.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# 4. Lower-case the "Unknown" placeholder text for the importance columns
#    (alumni, first-generation, residency, sat/act, gpa, volunteer work,
#    work experience). state.State (K2) keeps its original capitalized
#    "Unknown".
$ws.Range("D2:J2").Value = "unknown"
